# Append: 2025-12-26 01:24 JST
# A new scrape ran; one new listing ("航空会社とお客様のマッチングサービス...")
# sorts (by priority score, desc) into row 9, pushing the previous rows 9-13
# down to rows 10-14. Every row's "取得日時" (fetched-at) timestamp is bumped
# to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-26 01:24:48"

# 1) Make room: insert a blank row at row 9, shifting old rows 9-13 -> 10-14.
$ws.Rows.Item(9).Insert()

# 2) Populate the newly inserted row 9 with the new listing.
$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "【急募】航空会社とお客様のマッチングサービスのアプリ開発"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5461280"
$ws.Range("G9").Value = 85
$ws.Range("H9").Value = "◆開発 ◇アプリ"

# 3) Refresh the "取得日時" timestamp for every other data row (2-8 untouched
#    in content, 10-14 are the shifted former rows 9-13) to the new scrape time.
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A13").Value = $newTimestamp
$ws.Range("A14").Value = $newTimestamp

# 4) The row insert shifted cell *values* but left the worksheet's hyperlink
#    annotations pinned to their old cell addresses, so rebuild the whole
#    hyperlinks collection from scratch against the final layout.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5460562")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460357")
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460563")
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5460750")
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5460724")
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5460405")
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5460928")
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5461280")
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5460787")
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5016989")
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5460484")
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5461140")
$ws.Range("F13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5450323")
$ws.Range("F14").Style = "Hyperlink"
